$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number need an explicit text
# number format first, otherwise Excel auto-converts the assigned string
# into a numeric value (these Price cells are stored as text in the sheet).
$ws.Range("D2").Value = '70.421.75'
$ws.Range("E2").Value = '  -2.42%  '
$ws.Range("D3").Value = '2.554.87'
$ws.Range("E3").Value = '  -5.11%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.28'
$ws.Range("E5").Value = '  -3.66%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '170.11'
$ws.Range("E6").Value = '  -2.77%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.511'
$ws.Range("E8").Value = '  -2.40%  '
$ws.Range("D9").Value = '2.554.34'
$ws.Range("E9").Value = '  -5.12%  '
$ws.Range("E10").Value = '  -1.08%  '
$ws.Range("E11").Value = '  -0.09%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.345'
$ws.Range("E12").Value = '  -3.25%  '
$ws.Range("E13").Value = '  -3.40%  '
$ws.Range("D14").Value = '3.023.73'
$ws.Range("E14").Value = '  -5.04%  '
$ws.Range("E15").Value = '  -1.69%  '
$ws.Range("D16").Value = '70.315.18'
$ws.Range("E16").Value = '  -2.37%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.14'
$ws.Range("E17").Value = '  -4.32%  '
$ws.Range("D18").Value = '2.544.57'
$ws.Range("E18").Value = '  -4.96%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.69'
$ws.Range("E19").Value = '  -4.65%  '
$ws.Range("E20").Value = '  -5.74%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '361.56'
$ws.Range("E21").Value = '  -2.89%  '
$ws.Range("E22").Value = '  -4.95%  '
$ws.Range("E23").Value = '  +0.70%  '
$ws.Range("E24").Value = '  -0.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '69.87'
$ws.Range("E25").Value = '  -3.40%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.34'
$ws.Range("E27").Value = '  -4.73%  '
$ws.Range("E28").Value = '  -4.90%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.997'
$ws.Range("E29").Value = '  -0.30%  '
$ws.Range("D30").Value = '0.0₃0933'
$ws.Range("E30").Value = '  -4.59%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.89'
$ws.Range("E31").Value = '  -2.28%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '485.35'
$ws.Range("E32").Value = '  -3.60%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.29'
$ws.Range("E33").Value = '  -0.44%  '
$ws.Range("E34").Value = '  -3.17%  '
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("E36").Value = '  +6.72%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '157.74'
$ws.Range("E37").Value = '  -3.19%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.73'
$ws.Range("E38").Value = '  -4.82%  '
$ws.Range("E39").Value = '  -1.43%  '
$ws.Range("E40").Value = '  -3.77%  '
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("E42").Value = '  -5.09%  '
$ws.Range("E43").Value = '  -4.86%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.48'
$ws.Range("E44").Value = '  -3.57%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.322'
$ws.Range("E45").Value = '  -3.43%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '38.36'
$ws.Range("E46").Value = '  -2.90%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '145.40'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.57'
$ws.Range("E48").Value = '  -4.52%  '
$ws.Range("E49").Value = '  -5.71%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.596'
$ws.Range("E51").Value = '  -2.15%  '
